# Auto-generated script to apply price/profit updates to Seraph_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 613.95654
$arr[0,1] = 411.4
$arr[0,2] = 670.2222
$arr[0,3] = 1234.2
$arr[0,4] = 2010.6666
$arr[0,5] = -1066.2
$arr[0,6] = -2346.6666
$ws.Range("H17:N17").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2523.6428
$arr[0,1] = 1824
$arr[0,2] = 2640.25
$arr[0,3] = 1824
$arr[0,4] = 2640.25
$arr[0,5] = -1649
$arr[0,6] = -2990.25
$ws.Range("H19:N19").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 286.625
$arr[0,1] = 50.25
$ws.Range("H53:I53").Value = $arr
$ws.Range("K53").Value = 50.25
$ws.Range("M53").Value = 586.75

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3800
$arr[0,1] = 3000
$arr[0,2] = 4000
$arr[0,3] = 3000
$arr[0,4] = 4000
$arr[0,5] = -2685
$arr[0,6] = -4630
$ws.Range("H76:N76").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3800
$arr[0,1] = 3000
$arr[0,2] = 4000
$arr[0,3] = 3000
$arr[0,4] = 4000
$arr[0,5] = -1908
$arr[0,6] = -6184
$ws.Range("H79:N79").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 584.3333
$arr[0,1] = 515
$arr[0,2] = 639.8
$arr[0,3] = 1545
$arr[0,4] = 1919.4
$arr[0,5] = -547
$arr[0,6] = -3915.4
$ws.Range("H80:N80").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 584.3333
$arr[0,1] = 515
$arr[0,2] = 639.8
$arr[0,3] = 4635
$arr[0,4] = 5758.2
$arr[0,5] = 357
$arr[0,6] = -15742.2
$ws.Range("H83:N83").Value = $arr

$ws.Range("H86").Value = 2122.2
$ws.Range("J86").Value = 3502.5
$ws.Range("L86").Value = 3502.5
$ws.Range("N86").Value = -5748.5

$ws.Range("H89").Value = 2122.2
$ws.Range("J89").Value = 3502.5
$ws.Range("L89").Value = 17512.5
$ws.Range("N89").Value = -28744.5

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 100003510
$arr[0,1] = 100003510
$ws.Range("H107:I107").Value = $arr
$ws.Range("K107").Value = 100003510
$ws.Range("M107").Value = -100001590

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 1404.2572
$arr[0,1] = 1425.6818
$ws.Range("H137:I137").Value = $arr
$ws.Range("K137").Value = 4277.0454
$ws.Range("M137").Value = -1727.0454

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 6112.231
$arr[0,1] = 3311.5
$ws.Range("H32:I32").Value = $arr
$ws.Range("K32").Value = 3311.5
$ws.Range("M32").Value = -3024.5

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1773.5385
$arr[0,1] = 1764.48
$arr[0,2] = 2000
$arr[0,3] = 1764.48
$arr[0,4] = 2000
$arr[0,5] = -1552.48
$arr[0,6] = -2424
$ws.Range("H61:N61").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 598.2
$arr[0,1] = 524.4211
$ws.Range("H74:I74").Value = $arr
$ws.Range("K74").Value = 524.4211
$ws.Range("M74").Value = 349.5789

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 598.2
$arr[0,1] = 524.4211
$ws.Range("H77:I77").Value = $arr
$ws.Range("K77").Value = 2622.1055
$ws.Range("M77").Value = 1745.8945

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1773.5385
$arr[0,1] = 1764.48
$arr[0,2] = 2000
$arr[0,3] = 5293.440000000001
$arr[0,4] = 6000
$arr[0,5] = -2743.440000000001
$arr[0,6] = -11100
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2997.6667
$arr[0,1] = 3282.8572
$arr[0,2] = 1999.5
$arr[0,3] = 3282.8572
$arr[0,4] = 1999.5
$arr[0,5] = -3035.8572
$arr[0,6] = -2493.5
$ws.Range("H20:N20").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1741.4
$arr[0,1] = 1350
$arr[0,2] = 2002.3334
$arr[0,3] = 1350
$arr[0,4] = 2002.3334
$arr[0,5] = -1125
$arr[0,6] = -2452.3334
$ws.Range("H64:N64").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1741.4
$arr[0,1] = 1350
$arr[0,2] = 2002.3334
$arr[0,3] = 1350
$arr[0,4] = 2002.3334
$arr[0,5] = -570
$arr[0,6] = -3562.3334
$ws.Range("H67:N67").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 233.44444
$arr[0,1] = 189
$ws.Range("H7:I7").Value = $arr
$ws.Range("K7").Value = 189
$ws.Range("M7").Value = -76

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10214.5
$arr[0,1] = 14420
$arr[0,2] = 6009
$arr[0,3] = 14420
$arr[0,4] = 6009
$arr[0,5] = -14250
$arr[0,6] = -6349
$ws.Range("H15:N15").Value = $arr

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 2362.75
$arr[0,1] = 2343.1428
$ws.Range("H62:I62").Value = $arr
$ws.Range("K62").Value = 2343.1428
$ws.Range("M62").Value = -1719.1428

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 2362.75
$arr[0,1] = 2343.1428
$ws.Range("H65:I65").Value = $arr
$ws.Range("K65").Value = 11715.714
$ws.Range("M65").Value = -8595.714

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 12349.154
$arr[0,1] = 8774.462
$arr[0,2] = 15923.846
$arr[0,3] = 8774.462
$arr[0,4] = 15923.846
$arr[0,5] = -7276.462
$arr[0,6] = -18919.846
$ws.Range("H99:N99").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 771.2857
$arr[0,1] = 760.3
$ws.Range("H122:I122").Value = $arr
$ws.Range("K122").Value = 2280.9
$ws.Range("M122").Value = 169.1000000000004

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 12349.154
$arr[0,1] = 8774.462
$arr[0,2] = 15923.846
$arr[0,3] = 26323.386
$arr[0,4] = 47771.538
$arr[0,5] = -23853.386
$arr[0,6] = -52711.538
$ws.Range("H126:N126").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 71686.14
$arr[0,1] = 166715.33
$arr[0,2] = 414.25
$arr[0,3] = 1000291.98
$arr[0,4] = 2485.5
$arr[0,5] = -1000178.98
$arr[0,6] = -2711.5
$ws.Range("H2:N2").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 474.42856
$arr[0,1] = 474.42856
$ws.Range("H14:I14").Value = $arr
$ws.Range("K14").Value = 1423.28568
$ws.Range("M14").Value = -1250.28568

$ws.Range("H26").Value = 560.5
$ws.Range("J26").Value = 600.55554
$ws.Range("L26").Value = 1801.66662
$ws.Range("N26").Value = -2377.66662

$ws.Range("H132").Value = 3840.6155
$ws.Range("J132").Value = 8500
$ws.Range("L132").Value = 76500
$ws.Range("N132").Value = -81560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5179.4
$ws.Range("J80").Value = 5966.3335
$ws.Range("L80").Value = 5966.3335
$ws.Range("N80").Value = -7962.3335

$ws.Range("H83").Value = 5179.4
$ws.Range("J83").Value = 5966.3335
$ws.Range("L83").Value = 29831.6675
$ws.Range("N83").Value = -39815.6675

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2131.3845
$arr[0,1] = 1355.1666
$arr[0,2] = 2796.7144
$arr[0,3] = 1355.1666
$arr[0,4] = 2796.7144
$arr[0,5] = 564.8334
$arr[0,6] = -6636.7144
$ws.Range("H107:N107").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 74276.36
$arr[0,1] = 2359.875
$ws.Range("H122:I122").Value = $arr
$ws.Range("K122").Value = 7079.625
$ws.Range("M122").Value = -4629.625

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 7504.5
$arr[0,1] = 7504.5
$ws.Range("H4:I4").Value = $arr
$ws.Range("K4").Value = 7504.5
$ws.Range("M4").Value = -7391.5

$ws.Range("H18").Value = 49999
$ws.Range("J18").Value = 49999
$ws.Range("L18").Value = 49999
$ws.Range("N18").Value = -50343

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 7504.5
$arr[0,1] = 7504.5
$ws.Range("H28:I28").Value = $arr
$ws.Range("K28").Value = 7504.5
$ws.Range("M28").Value = -7272.5

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 7504.5
$arr[0,1] = 7504.5
$ws.Range("H37:I37").Value = $arr
$ws.Range("K37").Value = 7504.5
$ws.Range("M37").Value = -7397.5

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1498.3334
$arr[0,1] = 997
$arr[0,2] = 1749
$arr[0,3] = 997
$arr[0,4] = 1749
$arr[0,5] = -824
$arr[0,6] = -2095
$ws.Range("H55:N55").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2483.6667
$arr[0,1] = 1974.75
$arr[0,2] = 3501.5
$arr[0,3] = 1974.75
$arr[0,4] = 3501.5
$arr[0,5] = -1225.75
$arr[0,6] = -4999.5
$ws.Range("H68:N68").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2483.6667
$arr[0,1] = 1974.75
$arr[0,2] = 3501.5
$arr[0,3] = 9873.75
$arr[0,4] = 17507.5
$arr[0,5] = -6129.75
$arr[0,6] = -24995.5
$ws.Range("H71:N71").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1168.1538
$arr[0,1] = 818.7
$arr[0,2] = 2333
$arr[0,3] = 818.7
$arr[0,4] = 2333
$arr[0,5] = -457.7
$arr[0,6] = -3055
$ws.Range("H82:N82").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1168.1538
$arr[0,1] = 818.7
$arr[0,2] = 2333
$arr[0,3] = 818.7
$arr[0,4] = 2333
$arr[0,5] = 429.3
$arr[0,6] = -4829
$ws.Range("H85:N85").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 36166.5
$ws.Range("J64").Value = 36166.5
$ws.Range("L64").Value = 36166.5
$ws.Range("N64").Value = -36662.5

$ws.Range("H67").Value = 36166.5
$ws.Range("J67").Value = 36166.5
$ws.Range("L67").Value = 36166.5
$ws.Range("N67").Value = -37882.5

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 6364.615
$arr[0,1] = 967.5
$ws.Range("H81:I81").Value = $arr
$ws.Range("K81").Value = 1935
$ws.Range("M81").Value = -874

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 6364.615
$arr[0,1] = 967.5
$ws.Range("H84:I84").Value = $arr
$ws.Range("K84").Value = 9675
$ws.Range("M84").Value = -4371

$ws.Range("H113").Value = 570.8
$ws.Range("J113").Value = 544.2857
$ws.Range("L113").Value = 1632.8571
$ws.Range("N113").Value = -5972.8571

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 1294.6086
$arr[0,1] = 725.1053000000001
$ws.Range("H136:I136").Value = $arr
$ws.Range("K136").Value = 2175.3159
$ws.Range("M136").Value = 374.6840999999999
